$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C7"   = -12.526
    "A8"   = -22.188
    "A10"  = -21.736
    "A12"  = -21.683
    "C15"  = -13.636
    "A18"  = -21.784
    "C18"  = -10.732
    "C20"  = -12.46
    "C29"  = -12.124
    "C30"  = -13.347
    "C31"  = -13.358
    "A37"  = -20.029
    "C40"  = -12.782
    "C50"  = -13.326
    "A55"  = -21.868
    "A68"  = -21.736
    "C68"  = -11.001
    "C76"  = -13.045
    "A77"  = -20.843
    "A78"  = -19.951
    "A81"  = -21.869
    "A82"  = -22.152
    "C87"  = -13.416
    "C88"  = -12.985
    "C96"  = -12.894
    "C98"  = -13.23
    "C101" = -13.048
    "C102" = -13.091
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
